# This workbook drives the Excel "Solver" add-in to find portfolio weights
# (MPT!A13:A17) that minimize portfolio variance subject to:
#   - weights sum to 1                        (MPT!A18 = solver_lhs2)
#   - portfolio return equals a target value  (MPT!J19 = solver_lhs3, constrained
#     to equal MPT!K19 = solver_rhs3)
#
# The author changed the target return constraint (K19) from 350 to 100 and
# re-ran Solver (using the GRG Nonlinear engine, solver_eng=1), which produced
# a new optimal allocation in A13:A17. All other cells on the sheet recompute
# automatically from these two inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MPT")
$wsData = $wb.Worksheets.Item("Data")

# --- New Solver constraint: target annual return (solver_rhs3) ---
$ws.Range("K19").Value = 100

# --- New Solver result: optimal weights for BTC/ETH/XRP/BCH/LTC (A13:A17) ---
$ws.Range("A13").Value = 0.38254724025862341
$ws.Range("A14").Value = 0
$ws.Range("A15").Value = 0
$ws.Range("A16").Value = 0.61745275974137637
$ws.Range("A17").Value = 0

# --- Bump Solver bookkeeping defined names left behind by re-running Solver ---
$nEst = $ws.Names.Add("solver_est", "=1")
$nEst.Visible = $false

$nNwt = $ws.Names.Add("solver_nwt", "=1")
$nNwt.Visible = $false

$wb.Names.Item("MPT!solver_ver").RefersTo = "=3"

$excel.Calculate()

# --- Restore the on-screen selection/view left by the author after editing ---
$wsData.Activate() | Out-Null
$wsData.Range("G16").Select() | Out-Null

$ws.Activate() | Out-Null
$ws.Range("J19").Select() | Out-Null

$wb.Save()
